$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1208.909
$ws.Range("I28").Value = 1099.8823
$ws.Range("K28").Value = 1099.8823
$ws.Range("M28").Value = -614.8823
$ws.Range("H62").Value = 11025.929
$ws.Range("I62").Value = 7999.5
$ws.Range("J62").Value = 12236.5
$ws.Range("K62").Value = 7999.5
$ws.Range("L62").Value = 12236.5
$ws.Range("M62").Value = -7375.5
$ws.Range("N62").Value = -13484.5
$ws.Range("H65").Value = 11025.929
$ws.Range("I65").Value = 7999.5
$ws.Range("J65").Value = 12236.5
$ws.Range("K65").Value = 39997.5
$ws.Range("L65").Value = 61182.5
$ws.Range("M65").Value = -36877.5
$ws.Range("N65").Value = -67422.5
$ws.Range("H94").Value = 747.6667
$ws.Range("I94").Value = 747.6667
$ws.Range("K94").Value = 747.6667
$ws.Range("M94").Value = -296.6667
$ws.Range("H96").Value = 798.25
$ws.Range("I96").Value = 1242.1428
$ws.Range("J96").Value = 176.8
$ws.Range("K96").Value = 3726.4284
$ws.Range("L96").Value = 530.4000000000001
$ws.Range("M96").Value = -2353.4284
$ws.Range("N96").Value = -3276.4
$ws.Range("H99").Value = 708.5
$ws.Range("I99").Value = 708.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2125.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -627.5
$ws.Range("N99").ClearContents()
$ws.Range("H104").Value = 955.375
$ws.Range("I104").Value = 618.8333
$ws.Range("K104").Value = 1856.4999
$ws.Range("M104").Value = -109.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 69994
$ws.Range("J44").Value = 69994
$ws.Range("L44").Value = 69994
$ws.Range("N44").Value = -70970
$ws.Range("H55").Value = 50021
$ws.Range("J55").Value = 69994
$ws.Range("L55").Value = 69994
$ws.Range("N55").Value = -70624
$ws.Range("H63").Value = 3480.6667
$ws.Range("I63").Value = 1999
$ws.Range("J63").Value = 4962.3335
$ws.Range("K63").Value = 1999
$ws.Range("L63").Value = 4962.3335
$ws.Range("M63").Value = -1313
$ws.Range("N63").Value = -6334.3335
$ws.Range("H66").Value = 3480.6667
$ws.Range("I66").Value = 1999
$ws.Range("J66").Value = 4962.3335
$ws.Range("K66").Value = 9995
$ws.Range("L66").Value = 24811.6675
$ws.Range("M66").Value = -6563
$ws.Range("N66").Value = -31675.6675
$ws.Range("H110").Value = 3111.75
$ws.Range("I110").Value = 2998.75
$ws.Range("J110").Value = 3224.75
$ws.Range("K110").Value = 2998.75
$ws.Range("L110").Value = 3224.75
$ws.Range("M110").Value = -953.75
$ws.Range("N110").Value = -7314.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H128").Value = 5776.5557
$ws.Range("I128").Value = 5776.5557
$ws.Range("K128").Value = 17329.6671
$ws.Range("M128").Value = -14839.6671
$ws.Range("H134").Value = 1292.5294
$ws.Range("I134").Value = 1317.0625
$ws.Range("J134").Value = 900
$ws.Range("K134").Value = 3951.1875
$ws.Range("L134").Value = 2700
$ws.Range("M134").Value = -1416.1875
$ws.Range("N134").Value = -7770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 22.25
$ws.Range("I7").Value = 16.304348
$ws.Range("K7").Value = 16.304348
$ws.Range("M7").Value = 96.695652
$ws.Range("H22").Value = 479.125
$ws.Range("I22").Value = 479.125
$ws.Range("K22").Value = 479.125
$ws.Range("M22").Value = -129.125
$ws.Range("H122").Value = 3520.138
$ws.Range("I122").Value = 2811.7646
$ws.Range("K122").Value = 8435.293799999999
$ws.Range("M122").Value = -5985.293799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4191.5454
$ws.Range("J39").Value = 4254
$ws.Range("L39").Value = 12762
$ws.Range("N39").Value = -13350

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 24999
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 24999
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 24999
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -25489

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2817.8
$ws.Range("I40").Value = 2811.8572
$ws.Range("K40").Value = 2811.8572
$ws.Range("M40").Value = -2675.8572
$ws.Range("H56").Value = 11173.667
$ws.Range("I56").Value = 11173.667
$ws.Range("K56").Value = 11173.667
$ws.Range("M56").Value = -10482.667
$ws.Range("H68").Value = 2003.5
$ws.Range("I68").Value = 1815.4445
$ws.Range("K68").Value = 1815.4445
$ws.Range("M68").Value = -1066.4445
$ws.Range("H71").Value = 2003.5
$ws.Range("I71").Value = 1815.4445
$ws.Range("K71").Value = 9077.2225
$ws.Range("M71").Value = -5333.2225
$ws.Range("H93").Value = 47620284
$ws.Range("I93").Value = 62500884
$ws.Range("K93").Value = 62500884
$ws.Range("M93").Value = -62499636
$ws.Range("H122").Value = 16354.619
$ws.Range("I122").Value = 20503.6
$ws.Range("J122").Value = 12582.818
$ws.Range("K122").Value = 61510.8
$ws.Range("L122").Value = 37748.454
$ws.Range("M122").Value = -59060.8
$ws.Range("N122").Value = -42648.454
$ws.Range("H128").Value = 87995
$ws.Range("J128").Value = 87995
$ws.Range("L128").Value = 87995
$ws.Range("N128").Value = -97955

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 6305
$ws.Range("I18").Value = 6
$ws.Range("J18").Value = 12604
$ws.Range("K18").Value = 6
$ws.Range("L18").Value = 12604
$ws.Range("M18").Value = 167
$ws.Range("N18").Value = -12950
$ws.Range("H20").Value = 30011
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H62").Value = 5218.1
$ws.Range("J62").Value = 5710.125
$ws.Range("L62").Value = 5710.125
$ws.Range("N62").Value = -6958.125
$ws.Range("H65").Value = 5218.1
$ws.Range("J65").Value = 5710.125
$ws.Range("L65").Value = 28550.625
$ws.Range("N65").Value = -34790.625
$ws.Range("H122").Value = 6757.4443
$ws.Range("I122").Value = 6432.4165
$ws.Range("J122").Value = 7407.5
$ws.Range("K122").Value = 19297.2495
$ws.Range("L122").Value = 22222.5
$ws.Range("M122").Value = -16847.2495
$ws.Range("N122").Value = -27122.5
